# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> new value }
$updates = @{
    "展览" = @{
        3  = 1207
        4  = 14520
        5  = 17364
        8  = 59
        25 = 7137
        28 = 1160
        30 = 5836
        31 = 56
        33 = 132
        35 = 219
        36 = 5032
        37 = 31
    }
    "全部类型" = @{
        3  = 1207
        4  = 14520
        5  = 17364
        8  = 59
        26 = 7137
        29 = 1160
        32 = 5836
        33 = 56
        35 = 132
        37 = 219
        38 = 5032
        39 = 31
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
